$wb = $excel.ActiveWorkbook

# --- Update the data on the "Merge Comparison" sheet ---
$wsMerge = $wb.Worksheets.Item("Merge Comparison")

$row2 = @(17, 72, 84, 138, 75, 87, 127, 139, 124, 133)
$row3 = @(1040, 3460, 8200, 19460, 30366, 54235, 79193, 91954, 123620, 160781)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $wsMerge.Cells.Item(2, 2 + $i).Value = $row2[$i]
    $wsMerge.Cells.Item(3, 2 + $i).Value = $row3[$i]
}

# --- Move/resize the chart embedded on "Merge Comparison" ---
$chartObj = $wsMerge.ChartObjects(1)
$chartObj.Left = 182.3125
$chartObj.Top = 137
$chartObj.Width = 514.5
$chartObj.Height = 383

# --- Update the selected cell on "Algorithm Runtimes" (kept inactive) ---
$wsAlgo = $wb.Worksheets.Item("Algorithm Runtimes")
[void]$wsAlgo.Range("H5").Select()

# --- Update the selected cell on "Merge Comparison" and leave it as the active sheet/tab ---
[void]$wsMerge.Activate()
[void]$wsMerge.Range("M4").Select()
